# Update column F ("dSF") values per repull of data.
# Workbook has a single worksheet with a data table; the commit message
# indicates the data in this column was re-pulled / recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 2
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 7
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = -5
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = -1
$ws.Range("F21").Value = 8
$ws.Range("F22").Value = -6
$ws.Range("F24").Value = -2
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 2
$ws.Range("F28").Value = 1
$ws.Range("F29").Value = -4
$ws.Range("F30").Value = 1
$ws.Range("F31").Value = 1
$ws.Range("F32").Value = -2
$ws.Range("F33").Value = 1
$ws.Range("F34").Value = 2
$ws.Range("F35").Value = -1
$ws.Range("F36").Value = -2
$ws.Range("F37").Value = 2
$ws.Range("F39").Value = -1
